$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")

# Update confidential notice date from 2021-04-23 to 2021-04-26
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-26 for illustrative purposes only and are subject to change."
$ws.Rows.Item(42).AutoFit()

# Update performance/contribution figures in columns D and E
$ws.Range("D2").Value = 0.06299132906211981
$ws.Range("E2").Value = 0.002977963073258039
$ws.Range("D3").Value = 0.05696284327801913
$ws.Range("E3").Value = 0.00153168677005544
$ws.Range("D4").Value = 0.2904090734840309
$ws.Range("E4").Value = 0.004206098843322792
$ws.Range("D5").Value = 0.0364361523742425
$ws.Range("E5").Value = 0.02038983740810796
$ws.Range("D6").Value = 0.03260940099584093
$ws.Range("E6").Value = 0.001478260869565151
$ws.Range("D7").Value = 0.02948391533116265
$ws.Range("E7").Value = 0.002463546174845233
$ws.Range("D8").Value = 0.02888298623588812
$ws.Range("E8").Value = -0.008458192363460681
$ws.Range("D9").Value = 0.02441233551474594
$ws.Range("E9").Value = -0.01422444603288064
$ws.Range("D10").Value = 0.02508339118139279
$ws.Range("E10").Value = 0.004347958416125675
$ws.Range("D11").Value = 0.02298919145590069
$ws.Range("E11").Value = 0.006342775545445711
$ws.Range("D12").Value = 0.02221976227855735
$ws.Range("E12").Value = 0.005870342011230045
$ws.Range("D13").Value = 0.0219564877166712
$ws.Range("E13").Value = 0.006829854660692858
$ws.Range("D14").Value = 0.02191155436145285
$ws.Range("E14").Value = -0.02000895923547863
$ws.Range("D15").Value = 0.02119436566262793
$ws.Range("E15").Value = -0.01204112507332733
$ws.Range("D16").Value = 0.02182921289739981
$ws.Range("E16").Value = -0.01111638480177846
$ws.Range("D17").Value = 0.02110668018302708
$ws.Range("E17").Value = 0.001059267297059963
$ws.Range("D18").Value = 0.01550593376583505
$ws.Range("E18").Value = -0.008102633355840605
$ws.Range("D19").Value = 0.01660920032260404
$ws.Range("E19").Value = -0.001287001287001321
$ws.Range("D20").Value = 0.01562306586171978
$ws.Range("E20").Value = -0.006282722513089034
$ws.Range("D21").Value = 0.01575742968120766
$ws.Range("E21").Value = 0.001979485333813313
$ws.Range("D22").Value = 0.01590989771663303
$ws.Range("E22").Value = 0.01206471072114068
$ws.Range("D23").Value = 0.01544551367168222
$ws.Range("E23").Value = -0.01487057095648991
$ws.Range("D24").Value = 0.01472548937277071
$ws.Range("E24").Value = -0.01560509554140121
$ws.Range("D25").Value = 0.0141534615860711
$ws.Range("E25").Value = -0.005201309959545375
$ws.Range("D26").Value = 0.01479321659022361
$ws.Range("E26").Value = -0.005920038926283455
$ws.Range("D27").Value = 0.01277884991332334
$ws.Range("E27").Value = 0.01083032490974722
$ws.Range("D28").Value = 0.01329023948278654
$ws.Range("E28").Value = -0.0002954209748891046
$ws.Range("D29").Value = 0.01433548529932212
$ws.Range("E29").Value = 0.0005173305742369738
$ws.Range("D30").Value = 0.0130821500610472
$ws.Range("E30").Value = -0.006269174336401173
$ws.Range("D31").Value = 0.01272355571163474
$ws.Range("E31").Value = -0.01693753000068576
$ws.Range("D32").Value = 0.01336200197728575
$ws.Range("E32").Value = 0.00008978272580373314
$ws.Range("D33").Value = 0.01274056931215431
$ws.Range("E33").Value = -0.004494093477144134
$ws.Range("D34").Value = 0.006659406803368038
$ws.Range("E34").Value = 0.0139368827893418
$ws.Range("D35").Value = 0.005513606245300129
$ws.Range("E35").Value = 0.009395707645138973
$ws.Range("D36").Value = 0.005802728392591012
$ws.Range("E36").Value = 0.0213885651994139
$ws.Range("D37").Value = 0.005625830571804211
$ws.Range("E37").Value = -0.0002714019851116634
$ws.Range("D38").Value = 0.005083685647555631
$ws.Range("E38").Value = 0.01227125480016311
$ws.Range("E39").Value = 0.0008650761494950476

# Restore sheet protection to its original state
$ws.Protect("D382", $false, $true, $true, $true)
